$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update quantity (NÚMERO) from 1 to 6
$ws.Range("F6").Value = 6

# Row 9: first item line - update partida/description, entregado qty and unit price
$ws.Range("A9").Value = "2700"
$ws.Range("B9").Value = "QW"
$ws.Range("F9").Value = 2
$ws.Range("H9").Value = 6.428571428571429

# Row 10: second item line - cleared out entirely (data removed)
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("J10").Value = ""

# Row 26: update signature names - "RESPONSABLE DEL ALMACÉN" and "RECIBÍ DE CONFORMIDAD"
$ws.Range("F26").Value = "RE"
$ws.Range("J26").Value = "Daniel A. Benitez"
